# Fill in "NIL" bins for the days that have no lecture/practice data yet
# (rows 22-25, columns C "Practice" and D "Tech stack"), and move the
# active selection/view to reflect where the user was working afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 22..25) {
    $ws.Cells.Item($r, 3).Value = "NIL"
    $ws.Cells.Item($r, 4).Value = "NIL"
}

# Scroll the view down a bit and move the selection, matching where the
# user ended up after entering the new values.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("C26").Select()
